# Master Data F1.xlsx - "Perubahan database dan excel"
#
# Turns the raw attendance fractions in columns X (Hadir %) / Y (Tidak
# Hadir %) into their literal, Indonesian-locale-formatted percentage
# strings ("100%", "93,33%", ...), and fills in a previously-missing
# attendance mark (K11).
#
# Setting Range.Value directly to a string that *looks* like a percentage
# (e.g. "100%") makes Excel "smart" auto-convert it back into a number
# with a freshly minted percent NumberFormat/style - not the literal text
# cell the workbook actually needs. To store literal text while leaving
# each cell's existing style (border/alignment) untouched, we stage the
# text in a scratch cell formatted as Text ("@"), copy it, and
# PasteSpecial only the *values* into the destination - which carries
# over the text content but leaves the destination's own style alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("AZ1")
$scratch.NumberFormat = "@"

function Set-LiteralText($addr, $text) {
    $scratch.Value = $text
    $scratch.Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# Row 7 - Drs. Agus Setiawan, M.Kom / Aljabar Linier -> 100% / 0%
Set-LiteralText "X7" "100%"
Set-LiteralText "Y7" "0%"

# Row 8 - Novita Angra / Bahasa Indonesia untuk TI -> 100% / 0%
Set-LiteralText "X8" "100%"
Set-LiteralText "Y8" "0%"

# Row 9 - Fitria Nugrahani / Bahasa Inggris untuk TIK 2 -> 93,33% / 6,67%
Set-LiteralText "X9" "93,33%"
Set-LiteralText "Y9" "6,67%"

# Row 10 - Risna Sari / Basis Data 1 -> 100% / 0%
Set-LiteralText "X10" "100%"
Set-LiteralText "Y10" "0%"

# Row 11 - Herlino Nanang / Jaringan Komputer dan Komunikasi
# K11 was blank (no class held on P7) and now has an attendance mark.
$ws.Range("K11").Value = 1
Set-LiteralText "X11" "100%"
Set-LiteralText "Y11" "0%"

# Row 12 - Iklima Ermis Ismail / Pemrograman Web 1 -> 93,75% / 6,25%
Set-LiteralText "X12" "93,75%"
Set-LiteralText "Y12" "6,25%"

# Row 13 - Euis Oktavianti / Rekayasa Perangkat Lunak -> 100% / 0%
Set-LiteralText "X13" "100%"
Set-LiteralText "Y13" "0%"

# Row 14 - Dewi Yanti Liliana / Struktur Data -> 100% / 0%
Set-LiteralText "X14" "100%"
Set-LiteralText "Y14" "0%"

# Clean up the scratch cell so it doesn't linger in the saved sheet.
$scratch.Clear() | Out-Null
$excel.CutCopyMode = $false
